# no-op test
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
